$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.38"
$ws.Range("E2").Value = "'0.95%"
$ws.Range("D3").Value = "'29.89"
$ws.Range("E3").Value = "'12.35%"
$ws.Range("D4").Value = "'5.163"
$ws.Range("E4").Value = "'0.17%"
$ws.Range("D5").Value = "'0.05730"
$ws.Range("E5").Value = "'1.96%"
$ws.Range("D6").Value = "'6.583"
$ws.Range("E6").Value = "'1.38%"
$ws.Range("D7").Value = "'0.8564"
$ws.Range("E7").Value = "'4.64%"
$ws.Range("D8").Value = "'0.8733"
$ws.Range("E8").Value = "'4.96%"
$ws.Range("D9").Value = "'0.1360"
$ws.Range("E9").Value = "'2.48%"
$ws.Range("D10").Value = "'0.06988"
$ws.Range("E10").Value = "'0.95%"
$ws.Range("D11").Value = "'0.02905"
$ws.Range("E11").Value = "'0.27%"
$ws.Range("D12").Value = "'0.09376"
$ws.Range("E12").Value = "'-0.08%"
$ws.Range("D13").Value = "'0.001516"
$ws.Range("E13").Value = "'0.41%"
$ws.Range("D14").Value = "'0.04158"
$ws.Range("E14").Value = "'-9.23%"
$ws.Range("D15").Value = "'0.0006016"
$ws.Range("E15").Value = "'-93.96%"
$ws.Range("D16").Value = "'0.006042"
$ws.Range("E16").Value = "'-3.31%"
$ws.Range("D17").Value = "'3.510"
$ws.Range("E17").Value = "'-3.68%"
$ws.Range("D18").Value = "'3.030"
$ws.Range("E18").Value = "'0.17%"
$ws.Range("D19").Value = "'2.175"
$ws.Range("E19").Value = "'-0.74%"
$ws.Range("D20").Value = "'0.3143"
$ws.Range("E20").Value = "'1.01%"
$ws.Range("D21").Value = "'0.03318"
$ws.Range("E21").Value = "'6.76%"
$ws.Range("E22").Value = "'1.02%"
$ws.Range("D23").Value = "'3.616"
$ws.Range("E23").Value = "'-3.51%"
$ws.Range("E24").Value = "'2.64%"
$ws.Range("E25").Value = "'-1.39%"
$ws.Range("D26").Value = "'0.004498"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("D27").Value = "'0.0001178"
$ws.Range("E27").Value = "'20.22%"
$ws.Range("D28").Value = "'0.0001377"
$ws.Range("E28").Value = "'-1.59%"
$ws.Range("D40").Value = "'0.03787"
$ws.Range("E40").Value = "'4.15%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1068"
$ws.Range("E41").Value = "'-21.97%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002585"
$ws.Range("E42").Value = "'0.60%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003567"
$ws.Range("E43").Value = "'-42.08%"
$ws.Range("D44").Value = "'0.01001"
$ws.Range("E44").Value = "'23.66%"
$ws.Range("D45").Value = "'0.00005095"
$ws.Range("E45").Value = "'-4.92%"
$ws.Range("E46").Value = "'-0.23%"
$ws.Range("D47").Value = "'0.07981"
$ws.Range("E47").Value = "'-26.78%"
$ws.Range("D48").Value = "'0.002723"
$ws.Range("E48").Value = "'5.60%"
$ws.Range("D49").Value = "'0.00002095"
$ws.Range("E49").Value = "'-0.23%"
$ws.Range("E50").Value = "'-0.23%"
